$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 571
$ws1.Range("F3").Value = 189
$ws1.Range("F4").Value = 373
$ws1.Range("F5").Value = 422
$ws1.Range("F7").Value = 2415
$ws1.Range("F9").Value = 6307
$ws1.Range("F11").Value = 406

# Sheet "全部类型" (sheet4): update "想去人数" (F column) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 571
$ws4.Range("F3").Value = 189
$ws4.Range("F4").Value = 373
$ws4.Range("F5").Value = 422
$ws4.Range("F9").Value = 2415
$ws4.Range("F11").Value = 6307
$ws4.Range("F13").Value = 406
